$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A12").Value = "iii"
$ws.Range("B12").Value = "sdfsd@dsfsd.com"
